$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 7: North of Tyne Combined Authority (previously had all-zero/missing plan data)
$ws.Range("C7").Value = "North of Tyne Combined Authority"
$ws.Range("D7").Value = "NTCA"
$ws.Range("E7").Value = 0.3809523809523809
$ws.Range("F7").Value = 0.1666666666666667
$ws.Range("G7").Value = 0.2
$ws.Range("H7").Value = 0.3333333333333333
$ws.Range("I7").Value = 0.4
$ws.Range("J7").Value = 0.5
$ws.Range("K7").Value = 0.2
$ws.Range("L7").Value = 0.8
$ws.Range("M7").Value = 0
$ws.Range("N7").Value = 0.3271428571428571
$ws.Range("Q7").Value = "North East"

# Row 8: West of England Combined Authority (shifted down from row 7)
$ws.Range("C8").Value = "West of England Combined Authority"
$ws.Range("D8").Value = "WECA"
$ws.Range("E8").Value = 0.4761904761904762
$ws.Range("F8").Value = 0.3888888888888889
$ws.Range("G8").Value = 0.4
$ws.Range("H8").Value = 0.3333333333333333
$ws.Range("I8").Value = 0.2
$ws.Range("J8").Value = 0.5
$ws.Range("K8").Value = 0
$ws.Range("L8").Value = 0
$ws.Range("M8").Value = 0.5
$ws.Range("N8").Value = 0.3097619047619048
$ws.Range("Q8").Value = "South West"

# Row 9: Sheffield City Region (shifted down from row 8)
$ws.Range("C9").Value = "Sheffield City Region"
$ws.Range("D9").Value = "SCR"
$ws.Range("E9").Value = 0.2857142857142857
$ws.Range("F9").Value = 0.05555555555555555
$ws.Range("G9").Value = 0.2
$ws.Range("H9").Value = 0.5555555555555556
$ws.Range("I9").Value = 0.4
$ws.Range("J9").Value = 0.75
$ws.Range("K9").Value = 0
$ws.Range("L9").Value = 0.2
$ws.Range("M9").Value = 0
$ws.Range("N9").Value = 0.2620238095238095
$ws.Range("Q9").Value = "Yorkshire and The Humber"

# Row 10: Tees Valley Combined Authority (shifted down from row 9, values unchanged, all zero)
$ws.Range("C10").Value = "Tees Valley Combined Authority"
$ws.Range("D10").Value = "TVCA"
$ws.Range("Q10").Value = "North East"

# Row 11: Cambridgeshire and Peterborough Combined Authority (shifted down from row 10, values unchanged, all zero)
$ws.Range("C11").Value = "Cambridgeshire and Peterborough Combined Authority"
$ws.Range("D11").Value = "CPCA"
$ws.Range("Q11").Value = "East"

# Row 12: North East Combined Authority unchanged (still all zero)
